# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit: new functions now also fetch each team's season
# record and append it as three extra columns (AD:AF) on every table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cells (AD1:AF1) the same style as the rest of the
# header row (bold, bordered, centered) by copying the format from an
# existing header cell, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record for this team: 101 wins, 61 losses, 0 ties - same for
# every player row (2 through 45).
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 101
    $ws.Cells.Item($r, 31).Value = 61
    $ws.Cells.Item($r, 32).Value = 0
}
